# LocationToLocation.xlsx update
# - System Setup changes: clear the sample date value (44353) placed in column A
#   of the "Date" field on each tab, while keeping the date-format style where the
#   diff shows it retained, and clearing the cell completely on the first tab.
# - Update the saved selection/active cell on each tab to reflect selecting the
#   remaining data rows (mirrors what Excel records after selecting whole rows).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: RSTK-8171-Existing loc -----------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A5").Clear()
$ws1.Range("A7:XFD1048576").Select()

# --- Sheet 2: RSTK-8172-New loc -----------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A3").ClearContents()

# --- Sheet 3: RSTK-8173-SRL-Existing loc --------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").ClearContents()
$ws3.Range("A6:XFD1048576").Select()

# --- Sheet 4: RSTK-8174-SRL-New_loc -------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").ClearContents()
$ws4.Range("A6:XFD1048576").Select()
